$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.183.91"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = "'1.872.25"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.26%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'311.68"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('D8').Value = "'0.3901"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').Value = "'0.09524"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('E10').Value = '  +3.80%  '
$ws.Range('D11').Value = "'40.85"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').Value = "'6.446"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').Value = "'20.94"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.39%  '
$ws.Range('D14').Value = "'1.870.74"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.68%  '
$ws.Range('D15').Value = "'1.002"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = "'7.379"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = "'92.61"
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'0.06598"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('E20').Value = '  +3.09%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = "'6.168"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.89%  '
$ws.Range('D23').Value = "'28.242.68"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').Value = "'11.26"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('D25').Value = "'2.273"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('D26').Value = "'2.585"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.41%  '
$ws.Range('D27').Value = "'2.085.82"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.08%  '
$ws.Range('D28').Value = "'21.17"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.14%  '
$ws.Range('D29').Value = "'159.04"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.13%  '
$ws.Range('D30').Value = "'127.15"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').Value = "'0.1062"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('D32').Value = "'1.068"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.57%  '
$ws.Range('D33').Value = "'5.621"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('D34').Value = "'3.625"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').Value = "'0.06755"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('D36').Value = "'9.507"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.05%  '
$ws.Range('D37').Value = "'0.02412"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.25%  '
$ws.Range('D38').Value = "'0.2188"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.93%  '
$ws.Range('D39').Value = "'11.51"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('D40').Value = "'0.6359"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('E41').Value = '  +1.35%  '
$ws.Range('D42').Value = "'1.184"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.49%  '
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').Value = "'13.56"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.15%  '
$ws.Range('D45').Value = "'0.5982"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.07%  '
$ws.Range('D46').Value = "'1.278"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('D47').Value = "'3.659"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').Value = "'1.995"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.03%  '
$ws.Range('D49').Value = "'123.60"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('D51').Value = "'0.06851"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.88%  '
